$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Student"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "CTI"
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("C5").Value = "Student"

$ws.Range("J3").Select()
